$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pkm_move_targets")

# Update the four cells whose text changes a single apostrophe (')
# into a doubled apostrophe ('') as per the diff.
$ws.Range("D4").Value = "The user''s ally (if any)."
$ws.Range("B5").Value = "User''s field"
$ws.Range("D5").Value = "The user''s side of the field.  Affects the user and its ally (if any)."
$ws.Range("B7").Value = "Opponent''s field"
